# New Help pop-up updates
# Updates the "CreateSTP" sheet's mini-legend table (rows 1-21, cols A/F/G)
# with the new test-data labels, and resets the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateSTP")

# Preserve the shared-string insertion order used by the source workbook:
# Publications label, FullName/Short samples, community-organizer sample,
# then the Data1..Data6 header labels.

# Row 21: publications label update
$ws.Range("A21").Value = "Publications By Merck/EMD"

# Row 2 / Row 3: full-name & short-name sample values
$ws.Range("F2").Value = "FullName110012"
$ws.Range("G2").Value = "FullName210012"
$ws.Range("F3").Value = "Short110012"
$ws.Range("G3").Value = "Short210012"

# Rows 5/6: community organizer sample value
$ws.Range("F5").Value = "Uday "
$ws.Range("G5").Value = "Uday "
$ws.Range("G6").Value = "Uday "

# Row 1: header labels for the "Data#" scenario columns
$ws.Range("B1").Value = "Data1(Negative case)"
$ws.Range("C1").Value = "Data2(Negative Case)"
$ws.Range("D1").Value = "Data3 (Negative Case)"
$ws.Range("E1").Value = "Data4 (Negative Case)"
$ws.Range("F1").Value = "Data5(Mandatory fields)"
$ws.Range("G1").Value = "Data6 (All Fields)"

# Reset the sheet's selection to the header row
$ws.Range("A1:G1").Select() | Out-Null
